# edit.ps1 — khl_stats_1369_ext.xlsx update (2025-11-06 build)
#
# 1. Matches_SOG: append 4 newly played matches (rows 433-436).
# 2. Shots_HA / Shots_Summary: bump as_of_utc for every team, and refresh the
#    home/away shots-on-goal aggregates for the 8 teams that actually played
#    (Авангард, Адмирал, Амур, Барыс, Металлург Мг, Нефтехимик,
#     Салават Юлаев, Трактор).
# 3. Meta_ext: bump as_of_utc + build_version.

$wb = $excel.ActiveWorkbook

# ---- Matches_SOG: append 4 new matches (rows 433-436) ----
$wsMatches = $wb.Worksheets.Item("Matches_SOG")

$wsMatches.Cells.Item(433, 1).Value = "897731"
$wsMatches.Cells.Item(433, 2).Value = "2025-11-05T12:15:00"
$wsMatches.Cells.Item(433, 3).Value = "Амур"
$wsMatches.Cells.Item(433, 4).Value = "Авангард"
$wsMatches.Cells.Item(433, 5).Value = 41
$wsMatches.Cells.Item(433, 6).Value = 46
$wsMatches.Cells.Item(433, 7).Value = "khl_text"

$wsMatches.Cells.Item(434, 1).Value = "897732"
$wsMatches.Cells.Item(434, 2).Value = "2025-11-05T12:30:00"
$wsMatches.Cells.Item(434, 3).Value = "Адмирал"
$wsMatches.Cells.Item(434, 4).Value = "Металлург Мг"
$wsMatches.Cells.Item(434, 5).Value = 36
$wsMatches.Cells.Item(434, 6).Value = 20
$wsMatches.Cells.Item(434, 7).Value = "khl_text"

$wsMatches.Cells.Item(435, 1).Value = "897734"
$wsMatches.Cells.Item(435, 2).Value = "2025-11-05T17:00:00"
$wsMatches.Cells.Item(435, 3).Value = "Трактор"
$wsMatches.Cells.Item(435, 4).Value = "Нефтехимик"
$wsMatches.Cells.Item(435, 5).Value = 41
$wsMatches.Cells.Item(435, 6).Value = 19
$wsMatches.Cells.Item(435, 7).Value = "khl_text"

$wsMatches.Cells.Item(436, 1).Value = "897733"
$wsMatches.Cells.Item(436, 2).Value = "2025-11-05T17:30:00"
$wsMatches.Cells.Item(436, 3).Value = "Барыс"
$wsMatches.Cells.Item(436, 4).Value = "Салават Юлаев"
$wsMatches.Cells.Item(436, 5).Value = 35
$wsMatches.Cells.Item(436, 6).Value = 36
$wsMatches.Cells.Item(436, 7).Value = "khl_text"

# ---- Shots_HA: refresh as_of_utc for all teams + updated aggregates for the 8 teams that played ----
$wsHA = $wb.Worksheets.Item("Shots_HA")
for ($r = 2; $r -le 23; $r++) {
    $wsHA.Cells.Item($r, 4).Value = "2025-11-05T17:30:00Z"
}

# row 2
$wsHA.Cells.Item(2, 6).Value = 17
$wsHA.Cells.Item(2, 11).Value = 618
$wsHA.Cells.Item(2, 12).Value = 513
$wsHA.Cells.Item(2, 13).Value = 36.4
$wsHA.Cells.Item(2, 14).Value = 30.2

# row 4
$wsHA.Cells.Item(4, 5).Value = 15
$wsHA.Cells.Item(4, 7).Value = 580
$wsHA.Cells.Item(4, 8).Value = 412
$wsHA.Cells.Item(4, 9).Value = 38.7
$wsHA.Cells.Item(4, 10).Value = 27.5

# row 6
$wsHA.Cells.Item(6, 5).Value = 18
$wsHA.Cells.Item(6, 7).Value = 552
$wsHA.Cells.Item(6, 8).Value = 640
$wsHA.Cells.Item(6, 9).Value = 30.7
$wsHA.Cells.Item(6, 10).Value = 35.6

# row 7
$wsHA.Cells.Item(7, 5).Value = 28
$wsHA.Cells.Item(7, 7).Value = 885
$wsHA.Cells.Item(7, 8).Value = 895
$wsHA.Cells.Item(7, 9).Value = 31.6
$wsHA.Cells.Item(7, 10).Value = 32

# row 13
$wsHA.Cells.Item(13, 6).Value = 15
$wsHA.Cells.Item(13, 11).Value = 423
$wsHA.Cells.Item(13, 12).Value = 399
$wsHA.Cells.Item(13, 13).Value = 28.2
$wsHA.Cells.Item(13, 14).Value = 26.6

# row 14
$wsHA.Cells.Item(14, 6).Value = 19
$wsHA.Cells.Item(14, 11).Value = 517
$wsHA.Cells.Item(14, 12).Value = 722
$wsHA.Cells.Item(14, 13).Value = 27.2
$wsHA.Cells.Item(14, 14).Value = 38

# row 16
$wsHA.Cells.Item(16, 6).Value = 25
$wsHA.Cells.Item(16, 11).Value = 694
$wsHA.Cells.Item(16, 12).Value = 733
$wsHA.Cells.Item(16, 13).Value = 27.8
$wsHA.Cells.Item(16, 14).Value = 29.3

# row 21
$wsHA.Cells.Item(21, 5).Value = 18
$wsHA.Cells.Item(21, 7).Value = 599
$wsHA.Cells.Item(21, 8).Value = 541
$wsHA.Cells.Item(21, 9).Value = 33.3
$wsHA.Cells.Item(21, 10).Value = 30.1

# ---- Shots_Summary: refresh as_of_utc for all teams + updated aggregates for the 8 teams that played ----
$wsSummary = $wb.Worksheets.Item("Shots_Summary")
for ($r = 2; $r -le 23; $r++) {
    $wsSummary.Cells.Item($r, 4).Value = "2025-11-05T17:30:00Z"
}

# row 2
$wsSummary.Cells.Item(2, 5).Value = 38
$wsSummary.Cells.Item(2, 6).Value = 1305
$wsSummary.Cells.Item(2, 7).Value = 1103
$wsSummary.Cells.Item(2, 8).Value = 34.3
$wsSummary.Cells.Item(2, 9).Value = 29

# row 4
$wsSummary.Cells.Item(4, 5).Value = 35
$wsSummary.Cells.Item(4, 6).Value = 1218
$wsSummary.Cells.Item(4, 7).Value = 973
$wsSummary.Cells.Item(4, 9).Value = 27.8

# row 6
$wsSummary.Cells.Item(6, 5).Value = 39
$wsSummary.Cells.Item(6, 6).Value = 1144
$wsSummary.Cells.Item(6, 7).Value = 1409
$wsSummary.Cells.Item(6, 8).Value = 29.3
$wsSummary.Cells.Item(6, 9).Value = 36.1

# row 7
$wsSummary.Cells.Item(7, 5).Value = 43
$wsSummary.Cells.Item(7, 6).Value = 1300
$wsSummary.Cells.Item(7, 7).Value = 1410
$wsSummary.Cells.Item(7, 8).Value = 30.2
$wsSummary.Cells.Item(7, 9).Value = 32.8

# row 13
$wsSummary.Cells.Item(13, 5).Value = 40
$wsSummary.Cells.Item(13, 6).Value = 1307
$wsSummary.Cells.Item(13, 7).Value = 1045
$wsSummary.Cells.Item(13, 8).Value = 32.7
$wsSummary.Cells.Item(13, 9).Value = 26.1

# row 14
$wsSummary.Cells.Item(14, 5).Value = 42
$wsSummary.Cells.Item(14, 6).Value = 1241
$wsSummary.Cells.Item(14, 7).Value = 1492
$wsSummary.Cells.Item(14, 8).Value = 29.5
$wsSummary.Cells.Item(14, 9).Value = 35.5

# row 16
$wsSummary.Cells.Item(16, 5).Value = 40
$wsSummary.Cells.Item(16, 6).Value = 1097
$wsSummary.Cells.Item(16, 7).Value = 1151
$wsSummary.Cells.Item(16, 8).Value = 27.4
$wsSummary.Cells.Item(16, 9).Value = 28.8

# row 21
$wsSummary.Cells.Item(21, 5).Value = 42
$wsSummary.Cells.Item(21, 6).Value = 1425
$wsSummary.Cells.Item(21, 7).Value = 1325
$wsSummary.Cells.Item(21, 8).Value = 33.9
$wsSummary.Cells.Item(21, 9).Value = 31.5

# ---- Meta_ext: bump as_of_utc + build_version ----
$wsMeta = $wb.Worksheets.Item("Meta_ext")
$wsMeta.Cells.Item(2, 2).Value = "2025-11-05T17:30:00Z"
$wsMeta.Cells.Item(2, 4).Value = 46